# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 277
$wsExhibition.Range("F3").Value = 171
$wsExhibition.Range("F4").Value = 2056
$wsExhibition.Range("F5").Value = 1648
$wsExhibition.Range("F6").Value = 297
$wsExhibition.Range("F7").Value = 84
$wsExhibition.Range("F8").Value = 672

# --- Sheet "全部类型" (All Types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 277
$wsAll.Range("F3").Value = 171
$wsAll.Range("F4").Value = 2056
$wsAll.Range("F5").Value = 1648
$wsAll.Range("F6").Value = 297
$wsAll.Range("F8").Value = 84
$wsAll.Range("F9").Value = 672
